$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date (Overview!G2) and Correspond Handoff Datetime (de-de!H2)
# share the same underlying string "2016-11-23 11:58:12" -> "2016-11-23 12:00:04"
$wsOverview.Range("G2").Value = "2016-11-23 12:00:04"
$wsDeDe.Range("H2").Value = "2016-11-23 12:00:04"

# zh-cn Correspond Handoff Datetime: 2016-11-23 11:57:57 -> 2016-11-23 11:59:50
$wsZhCn.Range("H2").Value = "2016-11-23 11:59:50"

# zh-cn Correspond Handback DateTime: 2016-11-23 11:58:50 -> 2016-11-23 12:00:48
$wsZhCn.Range("K2").Value = "2016-11-23 12:00:48"

# de-de Correspond Handback DateTime: 2016-11-23 11:59:09 -> 2016-11-23 12:01:07
$wsDeDe.Range("K2").Value = "2016-11-23 12:01:07"
